$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "AMSIN" (sheet1): append rows 41 and 42 to the registration
# history table (dimension grows from A1:G40 to A1:G42).
# ------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Row 41
$wsAmsin.Range("A41").Value = "'2023-02-17"
$wsAmsin.Range("B41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAmsin.Range("B41").Value = 44974.45777782407
$wsAmsin.Range("C41").Value = "ocr173fstccycle"
$wsAmsin.Range("D41").Value = 42
$wsAmsin.Range("E41").Value = 41
$wsAmsin.Range("F41").Value = 1
$wsAmsin.Range("G41").Value = 1.35

# Row 42
$wsAmsin.Range("A42").Value = "'2023-02-20"
$wsAmsin.Range("B42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAmsin.Range("B42").Value = 44977.42915909722
$wsAmsin.Range("C42").Value = "173ocrflow"
$wsAmsin.Range("D42").Value = 42
$wsAmsin.Range("E42").Value = 41
$wsAmsin.Range("F42").Value = 1
$wsAmsin.Range("G42").Value = 1.31

# ------------------------------------------------------------------
# Sheet "AMS" (sheet3): append rows 29, 30 and 31 to the registration
# history table (dimension grows from A1:G28 to A1:G31).
# ------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 29
$wsAms.Range("A29").Value = "'2023-02-20"
$wsAms.Range("B29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Range("B29").Value = 44977.60965219908
$wsAms.Range("C29").Value = "173ocrflow"
$wsAms.Range("D29").Value = 42
$wsAms.Range("E29").Value = 41
$wsAms.Range("F29").Value = 1
$wsAms.Range("G29").Value = 1.32

# Row 30
$wsAms.Range("A30").Value = "'2023-02-20"
$wsAms.Range("B30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Range("B30").Value = 44977.82466577546
$wsAms.Range("C30").Value = "173liveocr"
$wsAms.Range("D30").Value = 42
$wsAms.Range("E30").Value = 41
$wsAms.Range("F30").Value = 1
$wsAms.Range("G30").Value = 1.2

# Row 31
$wsAms.Range("A31").Value = "'2023-02-21"
$wsAms.Range("B31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsAms.Range("B31").Value = 44978.4902851037
$wsAms.Range("C31").Value = "ocrrecheck173"
$wsAms.Range("D31").Value = 42
$wsAms.Range("E31").Value = 42
$wsAms.Range("F31").Value = 0
$wsAms.Range("G31").Value = 2.02
